$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text updates (shared strings with uniform run formatting, so a
#    plain .Value assignment reproduces the same visible text/format).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# ---------------------------------------------------------------------------
# 2. Cells that change TYPE (number <-> shared "N/A"-style text). We copy the
#    formatting AND value from an existing, untouched cell that already has
#    the exact desired style/type, then (for numeric destinations) overwrite
#    just the value.
# ---------------------------------------------------------------------------

# --- Cells that must become the text "0" (style 14, shared string "0") ---
foreach ($addr in @("C15", "D18", "C26", "F27")) {
    $ws.Range("D14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range("D14").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# --- Cell that must become the text "***.*" (style 14, shared string "***.*") ---
foreach ($addr in @("E18")) {
    $ws.Range("E14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range("E14").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# --- Cells that must become plain integers (style 16) ---
$numStyle16 = @{ "D15" = 2; "G15" = 2; "D26" = 2; "D27" = 1 }
foreach ($addr in $numStyle16.Keys) {
    $ws.Range("F15").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $numStyle16[$addr]
}

# --- Cells that must become percent-change numbers (style 15) ---
$numStyle15 = @{ "E15" = -100; "H15" = -50; "E26" = -100; "E27" = -100 }
foreach ($addr in $numStyle15.Keys) {
    $ws.Range("M15").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $numStyle15[$addr]
}

# ---------------------------------------------------------------------------
# 3. Plain value updates (style stays the same, only the number changes).
# ---------------------------------------------------------------------------
$values = @{
    "J15" = 4;       "K15" = -25;                 "N15" = -57.142857142857;
    "C16" = 2;       "E16" = 100;                 "F16" = 8;
    "G16" = 7;       "H16" = 14.285714285714;     "I16" = 31;
    "J16" = 27;      "K16" = 14.814814814814;     "L16" = 47.619047619047;
    "M16" = -24.390243902439; "N16" = -84.653465346534;
    "C17" = 2;       "D17" = 7;                   "E17" = -71.428571428571;
    "G17" = 16;      "H17" = -12.5;                "I17" = 48;
    "J17" = 42;      "K17" = 14.285714285714;     "L17" = 71.428571428571;
    "M17" = 77.777777777777;  "N17" = -32.394366197183;
    "G18" = 7;       "H18" = -14.285714285714;    "I18" = 26;
    "K18" = -42.222222222222; "L18" = 52.941176470588;
    "M18" = 13.043478260869;  "N18" = -84.337349397590;
    "C19" = 11;      "D19" = 8;                   "E19" = 37.5;
    "F19" = 27;      "G19" = 23;                  "H19" = 17.391304347826;
    "I19" = 104;     "J19" = 89;                  "K19" = 16.853932584269;
    "L19" = 121.276595744681; "M19" = 67.741935483871;
    "N19" = -42.541436464088;
    "C20" = 3;       "E20" = 200;                 "F20" = 7;
    "G20" = 4;       "H20" = 75;                  "I20" = 24;
    "J20" = 19;      "K20" = 26.315789473684;     "L20" = 500;
    "M20" = 700;     "N20" = -80.8;
    "C21" = 19;      "D21" = 19;                  "E21" = 0;
    "F21" = 63;      "G21" = 59;                  "H21" = 6.779661016949;
    "I21" = 236;     "J21" = 226;                 "K21" = 4.424778761061;
    "L21" = 101.709401709402; "M21" = 45.679012345679;
    "N21" = -68.824306472919;
    "C22" = 1;       "I22" = 12;                  "K22" = 100;
    "L22" = 50;      "M22" = 300;
    "C23" = 2;       "D23" = 4;                   "G23" = 15;
    "H23" = -53.333333333333; "I23" = 38;         "J23" = 33;
    "K23" = 15.151515151515;  "L23" = 100;        "M23" = 46.153846153846;
    "C24" = 9;       "D24" = 13;                  "E24" = -30.769230769230;
    "F24" = 39;      "G24" = 41;                  "H24" = -4.878048780487;
    "I24" = 128;     "J24" = 138;                 "K24" = -7.246376811594;
    "L24" = 36.170212765957;  "M24" = 10.344827586206;
    "C25" = 6;       "D25" = 1;                   "E25" = 500;
    "F25" = 27;      "G25" = 16;                  "H25" = 68.75;
    "I25" = 68;      "J25" = 65;                  "K25" = 4.615384615384;
    "L25" = 100;     "M25" = 15.254237288135;
    "G26" = 3;       "H26" = -66.666666666666;    "J26" = 7;
    "K26" = -42.857142857142;
    "H27" = -100;    "J27" = 7;                   "K27" = 57.142857142857;
    "L27" = 57.142857142857;
    "M28" = -50;
    "M29" = -50;
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
